# Auto-generated Excel COM-interop script
# Applies updated market-price/profit figures scraped by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 8149.8125
$ws.Range("I2").Value = 1643.75
$ws.Range("J2").Value = 14655.875
$ws.Range("K2").Value = 1643.75
$ws.Range("L2").Value = 14655.875
$ws.Range("M2").Value = -1530.75
$ws.Range("N2").Value = -14881.875
$ws.Range("H12").Value = 632
$ws.Range("J12").Value = 450
$ws.Range("L12").Value = 450
$ws.Range("N12").Value = -790
$ws.Range("H17").Value = 1880
$ws.Range("J17").Value = 1981.875
$ws.Range("L17").Value = 5945.625
$ws.Range("N17").Value = -6281.625
$ws.Range("H28").Value = 9012.462
$ws.Range("I28").Value = 4045.375
$ws.Range("J28").Value = 16959.8
$ws.Range("K28").Value = 4045.375
$ws.Range("L28").Value = 16959.8
$ws.Range("M28").Value = -3560.375
$ws.Range("N28").Value = -17929.8
$ws.Range("H43").Value = 15873.25
$ws.Range("I43").Value = 25996
$ws.Range("K43").Value = 25996
$ws.Range("M43").Value = -25927
$ws.Range("H86").Value = 2599.818
$ws.Range("I86").Value = 2650
$ws.Range("J86").Value = 2539.6
$ws.Range("K86").Value = 2650
$ws.Range("L86").Value = 2539.6
$ws.Range("M86").Value = -1527
$ws.Range("N86").Value = -4785.6
$ws.Range("H89").Value = 2599.818
$ws.Range("I89").Value = 2650
$ws.Range("J89").Value = 2539.6
$ws.Range("K89").Value = 13250
$ws.Range("L89").Value = 12698
$ws.Range("M89").Value = -7634
$ws.Range("N89").Value = -23930
$ws.Range("H95").Value = 29642.25
$ws.Range("J95").Value = 29642.25
$ws.Range("L95").Value = 29642.25
$ws.Range("N95").Value = -35134.25
$ws.Range("H96").Value = 7149525
$ws.Range("I96").Value = 8614.143
$ws.Range("J96").Value = 14290435
$ws.Range("K96").Value = 25842.429
$ws.Range("L96").Value = 42871305
$ws.Range("M96").Value = -24469.429
$ws.Range("N96").Value = -42874051
$ws.Range("H99").Value = 328.6
$ws.Range("I99").Value = 324.5
$ws.Range("K99").Value = 973.5
$ws.Range("M99").Value = 524.5
$ws.Range("H106").Value = 3116
$ws.Range("I106").Value = 1661.7
$ws.Range("K106").Value = 1661.7
$ws.Range("M106").Value = -1030.7
$ws.Range("H112").Value = 2842.739
$ws.Range("I112").Value = 1272.25
$ws.Range("J112").Value = 3173.3684
$ws.Range("K112").Value = 3816.75
$ws.Range("L112").Value = 9520.1052
$ws.Range("M112").Value = -2708.75
$ws.Range("N112").Value = -11736.1052
$ws.Range("H121").Value = 1105
$ws.Range("J121").Value = 1105
$ws.Range("L121").Value = 3315
$ws.Range("N121").Value = -6809
$ws.Range("H140").Value = 120000
$ws.Range("J140").Value = 120000
$ws.Range("L140").Value = 120000
$ws.Range("N140").Value = -130360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11502133
$ws.Range("I2").Value = 13144867
$ws.Range("K2").Value = 13144867
$ws.Range("M2").Value = -13144754
$ws.Range("H32").Value = 54494.953
$ws.Range("I32").Value = 54494.953
$ws.Range("K32").Value = 54494.953
$ws.Range("M32").Value = -54207.953
$ws.Range("H45").Value = 2518.4546
$ws.Range("I45").Value = 2267.111
$ws.Range("K45").Value = 2267.111
$ws.Range("M45").Value = -1890.111
$ws.Range("H61").Value = 4145.9546
$ws.Range("I61").Value = 2980.125
$ws.Range("K61").Value = 2980.125
$ws.Range("M61").Value = -2768.125
$ws.Range("H74").Value = 245782.1
$ws.Range("I74").Value = 526804.0600000001
$ws.Range("J74").Value = 3081.318
$ws.Range("K74").Value = 526804.0600000001
$ws.Range("L74").Value = 3081.318
$ws.Range("M74").Value = -525930.0600000001
$ws.Range("N74").Value = -4829.318
$ws.Range("H77").Value = 245782.1
$ws.Range("I77").Value = 526804.0600000001
$ws.Range("J77").Value = 3081.318
$ws.Range("K77").Value = 2634020.3
$ws.Range("L77").Value = 15406.59
$ws.Range("M77").Value = -2629652.3
$ws.Range("N77").Value = -24142.59
$ws.Range("H97").Value = 1484162.9
$ws.Range("I97").Value = 1613170.8
$ws.Range("J97").Value = 572.5
$ws.Range("K97").Value = 1613170.8
$ws.Range("L97").Value = 572.5
$ws.Range("M97").Value = -1612674.8
$ws.Range("N97").Value = -1564.5
$ws.Range("H116").Value = 11502133
$ws.Range("I116").Value = 13144867
$ws.Range("K116").Value = 13144867
$ws.Range("M116").Value = -13142573
$ws.Range("H136").Value = 4145.9546
$ws.Range("I136").Value = 2980.125
$ws.Range("K136").Value = 8940.375
$ws.Range("M136").Value = -6390.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11502133
$ws.Range("I3").Value = 13144867
$ws.Range("K3").Value = 13144867
$ws.Range("M3").Value = -13144753
$ws.Range("H94").Value = 1114.2703
$ws.Range("I94").Value = 854.1613
$ws.Range("K94").Value = 854.1613
$ws.Range("M94").Value = -403.1613
$ws.Range("H105").Value = 58839396
$ws.Range("I105").Value = 58839396
$ws.Range("K105").Value = 58839396
$ws.Range("M105").Value = -58837649

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 101053.5
$ws.Range("I122").Value = 126067
$ws.Range("K122").Value = 378201
$ws.Range("M122").Value = -375751
$ws.Range("H134").Value = 10986.529
$ws.Range("I134").Value = 11447.846
$ws.Range("J134").Value = 9487.25
$ws.Range("K134").Value = 34343.538
$ws.Range("L134").Value = 28461.75
$ws.Range("M134").Value = -31808.538
$ws.Range("N134").Value = -33531.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 29692982
$ws.Range("I4").Value = 38517550
$ws.Range("K4").Value = 115552650
$ws.Range("M4").Value = -115552538
$ws.Range("H11").Value = 8230.6
$ws.Range("I11").Value = 9963.166999999999
$ws.Range("K11").Value = 29889.501
$ws.Range("M11").Value = -29749.501
$ws.Range("H109").Value = 3361.25
$ws.Range("I109").Value = 2848.6365
$ws.Range("J109").Value = 9000
$ws.Range("K109").Value = 8545.9095
$ws.Range("L109").Value = 27000
$ws.Range("M109").Value = -7505.9095
$ws.Range("N109").Value = -29080
$ws.Range("H134").Value = 113826.555
$ws.Range("I134").Value = 167223.33
$ws.Range("K134").Value = 501669.99
$ws.Range("M134").Value = -496599.99

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 10000
$ws.Range("J107").Value = 10000
$ws.Range("L107").Value = 10000
$ws.Range("N107").Value = -13840
$ws.Range("H126").Value = 2570.6
$ws.Range("I126").Value = 2309.75
$ws.Range("K126").Value = 6929.25
$ws.Range("M126").Value = -4459.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8996.218999999999
$ws.Range("I7").Value = 7663.5
$ws.Range("K7").Value = 7663.5
$ws.Range("M7").Value = -7551.5
$ws.Range("H40").Value = 27789718
$ws.Range("I40").Value = 27789718
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 27789718
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -27789582
$ws.Range("N40").ClearContents()
$ws.Range("H61").Value = 2261.2222
$ws.Range("I61").Value = 923.5
$ws.Range("K61").Value = 923.5
$ws.Range("M61").Value = -721.5
$ws.Range("H68").Value = 2011.125
$ws.Range("I68").Value = 1777.8
$ws.Range("K68").Value = 1777.8
$ws.Range("M68").Value = -1028.8
$ws.Range("H71").Value = 2011.125
$ws.Range("I71").Value = 1777.8
$ws.Range("K71").Value = 8889
$ws.Range("M71").Value = -5145
$ws.Range("H113").Value = 2261.2222
$ws.Range("I113").Value = 923.5
$ws.Range("K113").Value = 923.5
$ws.Range("M113").Value = 1246.5
$ws.Range("H122").Value = 4372.8
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 8996.218999999999
$ws.Range("I126").Value = 7663.5
$ws.Range("K126").Value = 22990.5
$ws.Range("M126").Value = -20520.5
$ws.Range("H132").Value = 4991.6924
$ws.Range("I132").Value = 3995.95
$ws.Range("K132").Value = 11987.85
$ws.Range("M132").Value = -9457.849999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 831.05554
$ws.Range("J113").Value = 1147.2858
$ws.Range("L113").Value = 3441.8574
$ws.Range("N113").Value = -7781.857400000001
$ws.Range("H114").Value = 48199
$ws.Range("J114").Value = 48199
$ws.Range("L114").Value = 48199
$ws.Range("N114").Value = -56877
$ws.Range("H119").Value = 20100
$ws.Range("J119").Value = 20100
$ws.Range("L119").Value = 20100
$ws.Range("N119").Value = -29776
$ws.Range("H122").Value = 6313.3335
$ws.Range("I122").Value = 6313.3335
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 18940.0005
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -16490.0005
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 3719.1538
$ws.Range("I126").Value = 2185.25
$ws.Range("J126").Value = 8832.166999999999
$ws.Range("K126").Value = 6555.75
$ws.Range("L126").Value = 26496.501
$ws.Range("M126").Value = -4085.75
$ws.Range("N126").Value = -31436.501

